$wb = $excel.ActiveWorkbook

# D2/D3 Bleu Score values for each of the 8 worksheets, in sheet order.
$values = @(
    @(0.412073170027898,  0.3926788350602298),
    @(0.3209398737453965, 0.2884115829516117),
    @(0.412073170027898,  0.3948788672519359),
    @(0.412073170027898,  0.3948788672519359),
    @(0.2175752258356586, 0.2042477121329563),
    @(0.2175752258356586, 0.1970724348722688),
    @(0.2255626915519154, 0.1901600827652285),
    @(0.2255626915519154, 0.2114603883333068)
)

for ($i = 0; $i -lt $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i + 1)

    # New header cell, copying the same (bold/border/centered) style used
    # by the existing "Cosine Similarity"/"Euclidean Distance" headers.
    $ws.Range("D1").Value = "Bleu Score"
    $ws.Range("C1").Copy()
    $ws.Range("D1").PasteSpecial(-4122)

    # New data column.
    $pair = $values[$i]
    $ws.Range("D2").Value = $pair[0]
    $ws.Range("D3").Value = $pair[1]
}

$excel.CutCopyMode = $false
